# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to match the refreshed data snapshot (gh-pages output regenerated
# at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 25
$wsExpo.Range("F3").Value = 1001
$wsExpo.Range("F4").Value = 241
$wsExpo.Range("F7").Value = 951
$wsExpo.Range("F8").Value = 287
$wsExpo.Range("F22").Value = 761
$wsExpo.Range("F28").Value = 3342

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 25
$wsAll.Range("F6").Value = 1001
$wsAll.Range("F7").Value = 241
$wsAll.Range("F11").Value = 951
$wsAll.Range("F12").Value = 287
$wsAll.Range("F34").Value = 761
$wsAll.Range("F42").Value = 3342
